# Rebuild the "Out of PO" player table: remove the "Naji Marshall" row and
# re-sort the remaining rows into their new order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Kessler Edwards", "SF,PF", "Dallas Mavericks"),
    @("Jaden McDaniels", "SF,PF", "Minnesota Timberwolves"),
    @("Domantas Sabonis", "C", "Sacramento Kings"),
    @("Alperen Sengün", "C", "Houston Rockets"),
    @("Donovan Mitchell", "PG,SG", "Cleveland Cavaliers"),
    @("Malik Beasley", "SG,SF", "Detroit Pistons"),
    @("Josh Hart", "SG,SF,PF", "New York Knicks"),
    @("Dyson Daniels", "PG,SG,SF", "Atlanta Hawks"),
    @("Bilal Coulibaly", "SG,SF", "Washington Wizards"),
    @("Kristaps Porzingis", "PF,C", "Boston Celtics"),
    @("De'Andre Hunter", "SF,PF", "Cleveland Cavaliers"),
    @("Victor Wembanyama", "C", "San Antonio Spurs"),
    @("Kelly Oubre Jr.", "SG,SF", "Philadelphia 76ers"),
    @("Julius Randle", "PF,C", "Minnesota Timberwolves"),
    @("Cam Thomas", "SG,SF", "Brooklyn Nets"),
    @("Michael Porter Jr.", "SF,PF", "Denver Nuggets")
)

# Clear the old data range (rows 2-18) before writing the new, shorter table.
$ws.Range("A2:C18").ClearContents()

# Write column-by-column (all names, then all positions, then all teams) so
# freshly-introduced shared strings land in the same first-use order Excel
# itself would produce.
for ($i = 0; $i -lt $data.Count; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $data[$i][0]
}
for ($i = 0; $i -lt $data.Count; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $data[$i][1]
}
for ($i = 0; $i -lt $data.Count; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $data[$i][2]
}
